$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the histologySubgroup (O2) and subsite (S2) values for the second row
$ws.Range("O2").Value = ""
$ws.Range("S2").Value = ""
